# Map city stats & sms
# Adds a new SMS campaign row (row 9) to Sheet1:
#   - B9: batch date (serial 42074 = 2015-03-11), same date style as B2:B8
#   - C9: new SMS text (becomes a new shared string)
#   - D9: sql id (23934)
#   - E9: =+LEN(C9) formula (same pattern as E2:E8), evaluates to 160
# Also clears the leftover placeholder formatting that used to sit on the
# empty D9/E9 cells, and moves the saved selection to D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy B8's cell (value + its date number-format style) down into B9 so the
# new date reuses the existing style instead of creating a new one, then
# overwrite the value with the new batch date (3/11/2015 -> serial 42074).
$ws.Range("B8").Copy($ws.Range("B9")) | Out-Null
$ws.Range("B9").Value = 42074

# New SMS body text for row 9 (added to sharedStrings.xml on write).
$ws.Range("C9").Value = "ΠΡΟΛΑΒΕ ΤΗΝ 1+1 ΠΡΟΣΦΟΡΑ:ΜΕΧΡΙ ΤΗΝ ΚΥΡΙΑΚΗ 15/3 ΜΕ ΚΑΘΕ Goody's Extreme Burger Η CLICKDELIVERY ΣΟΥ ΚΑΝΕΙ ΔΩΡΟ ΑΛΛΟ ΕΝΑ! ΠΑΡΑΓΓΕΙΛΕ ΤΩΡΑ ΣΤΟ www.clickdelivery.gr"

# D9/E9 used to be empty placeholder cells carrying a leftover numeric
# style; clear that formatting before writing the real data so the cells
# end up on the default style like the rest of the sheet's D/E columns.
$ws.Range("D9:E9").ClearFormats()
$ws.Range("D9").Value = 23934
$ws.Range("E9").Formula = "=+LEN(C9)"

# Recalculate so the new formula carries a cached value.
$wb.Application.Calculate()

# Match the saved cursor position recorded in the sheet view.
$ws.Range("D9").Select()
